$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column header in H1 (matching dimension A1:G2 -> A1:H2)
$ws.Range("H1").Value = "Save"

# Copy the header formatting/style (bold, border, centered) from the
# neighboring "sum" header cell (G1) onto the new "Save" header cell (H1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the Save value for the existing data row
$ws.Range("H2").Value = 1
